$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The ko_KR column (D) for these two story lines had been mistakenly filled
# in with the zh_CN shared strings. Add the correct Korean translations as
# new shared strings and re-point D2/D3 at them.

$korean1 = @"
지마는 학생자치단의 또 다른 멤버 나탈리야를 찾아갔다. 서로에게 묻고 싶은 게 있는 두 사람… 그녀들은 둘 다 마음에 두고 있는 그 일을 속 시원히 털어놓을 수 있는 사람을 찾고 있었다.

"@

$korean2 = @"
문제는 아직 해결되지 않았지만, 그녀들은 이후 함께 걸어나갈 것이다.

"@

$ws.Range("D2").Value = $korean1
$ws.Range("D3").Value = $korean2

# Re-run row autofit so the newly entered multi-line text doesn't leave a
# stray custom row height behind (matches the source edit, which only
# touched the shared strings / cell references).
$ws.Rows(2).AutoFit()
$ws.Rows(3).AutoFit()
